# Second commit - some corrections and updates:
# add the "Days" header label to cell A1 of the timetable sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Select() | Out-Null
$ws.Range("A1").Value = "Days"
